# AutoCommit_8 мая 2024 г. 9:08:35_SibNout2023
# Fill in homework ("ДЗ") scores of 2 for a number of students/assignments
# that were previously left blank, and move the frozen-pane view back to
# the top of the gradebook (C4) instead of scrolled down to C7/G18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that receive a new score of 2 (previously empty).
$cellsToSet = @(
    "F5",
    "C6", "D6", "E6", "F6",
    "C7", "D7", "E7", "F7",
    "C8", "D8", "E8", "F8",
    "E9", "F9",
    "C10", "D10", "E10", "F10",
    "F12",
    "C14", "D14", "E14", "F14",
    "C15", "D15", "E15", "F15",
    "C17", "D17", "E17", "F17",
    "C21", "D21", "E21", "F21",
    "E23", "F23",
    "D24", "E24", "F24",
    "C25", "D25", "E25", "F25",
    "C26", "D26", "E26", "F26",
    "C28", "D28", "E28", "F28",
    "C29", "D29", "E29", "F29",
    "C30", "D30", "E30", "F30",
    "F31",
    "E32"
)

foreach ($addr in $cellsToSet) {
    $ws.Range($addr).Value = 2
}

# Restore the view to the top-left of the scrollable area and select C4,
# matching the saved workbook view (frozen pane at C3/row3-col2, scrolled
# to C4 with C4 as the active cell in the bottom-right pane).
$ws.Range("C4").Select()
